$wb = $excel.ActiveWorkbook

# --- Add a new sheet "2024-05-24" after the last existing sheet ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2024-05-24"

# --- Append a new row of data to the "current" sheet ---
$ws = $wb.Worksheets.Item("current")

# Enter the date as a formula that evaluates to text, so Excel doesn't
# auto-convert the "2024-05-24" string into a date serial number/style.
# Converting to a static value afterwards keeps the cell a plain string
# (matching the other date-like cells A2/A3 on this sheet) with no
# extra number-format style attached.
$ws.Range("A4").Formula = '="2024-05-24"'
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163) # xlPasteValues

$ws.Range("B4").Value = "Vrachhhh"
$ws.Range("C4").Value = 0
